# GBDS UPDATED DOS - November: rename sheet date from 25-11-2025 to 08-11-2025
# and fill in the daily delivery figures, updating the scroll/selection state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename the sheet (this also renames the "label" side of the
#     Print_Area defined name automatically, but not the RefersTo target,
#     so we fix that explicitly below).
$ws.Name = "08-11-2025"

foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "='08-11-2025'!`$A`$1:`$V`$97"
    }
}

# --- Enter the delivered-quantity figures for the day ---
$ws.Range("F12").Value = 72
$ws.Range("M13").Value = 324
$ws.Range("F14").Value = 63
$ws.Range("F16").Value = 1512
$ws.Range("M16").Value = 1836
$ws.Range("F19").Value = 63
$ws.Range("P77").Value = 3420

# --- Update the view: scroll down to row 65 and select P77:Q77 ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 65
$ws.Range("P77:Q77").Select()
